$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 391, shifting existing data (rows 391-396) down to 393-398
$ws.Range("A391:A392").EntireRow.Insert()

# New row 391: Bola 8 variety, Region del Maule
$ws.Range("A391").Value = 10
$ws.Range("B391").Value = "Vega Modelo de Temuco"
$ws.Range("C391").Value = "La Araucanía"
$ws.Range("D391").Value = 44656
$ws.Range("D391").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E391").Value = 9
$ws.Range("F391").Value = 100112032
$ws.Range("G391").Value = "Zapallo italiano"
$ws.Range("H391").Value = "Bola 8"
$ws.Range("I391").Value = "Primera"
$ws.Range("J391").Value = 65
$ws.Range("K391").Value = 10000
$ws.Range("L391").Value = 10000
$ws.Range("M391").Value = 10000
$ws.Range("N391").Value = "$/caja 60 unidades"
$ws.Range("O391").Value = "Región del Maule"
$ws.Range("P391").Value = 167
$ws.Range("Q391").Value = 60
$ws.Range("R391").Value = "Hortaliza"

# New row 392: Sin especificar, Region del Maule
$ws.Range("A392").Value = 10
$ws.Range("B392").Value = "Vega Modelo de Temuco"
$ws.Range("C392").Value = "La Araucanía"
$ws.Range("D392").Value = 44656
$ws.Range("D392").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E392").Value = 9
$ws.Range("F392").Value = 100112032
$ws.Range("G392").Value = "Zapallo italiano"
$ws.Range("H392").Value = "Sin especificar"
$ws.Range("I392").Value = "Primera"
$ws.Range("J392").Value = 150
$ws.Range("K392").Value = 10000
$ws.Range("L392").Value = 10000
$ws.Range("M392").Value = 10000
$ws.Range("N392").Value = "$/caja 60 unidades"
$ws.Range("O392").Value = "Región del Maule"
$ws.Range("P392").Value = 167
$ws.Range("Q392").Value = 60
$ws.Range("R392").Value = "Hortaliza"
